$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value2 = "Datos actualizados a 23 de Marzo de 2020 a las 08:46"

# Swap the province rows 47 (Lugo) and 48 (Almeria), keeping their
# "Casos activos" (column C) values attached to the correct province.
$nameA47 = $ws.Range("A47").Value2
$nameA48 = $ws.Range("A48").Value2
$ws.Range("A47").Value2 = $nameA48
$ws.Range("A48").Value2 = $nameA47

$valC47 = $ws.Range("C47").Value2
$valC48 = $ws.Range("C48").Value2
$ws.Range("C47").Value2 = $valC48
$ws.Range("C48").Value2 = $valC47
